# Auto-generated edit script: add data for 2024-12-28
# Updates 2024 (column K) year-to-date totals (and a few 2022/2023 column
# corrections) across the Citywide Totals, By Neighborhood, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 7836
$ws.Range("I3").Value = 7492
$ws.Range("K3").Value = 8117
$ws.Range("I4").Value = 1817
$ws.Range("J4").Value = 1850
$ws.Range("K4").Value = 1708
$ws.Range("K5").Value = 578
$ws.Range("K6").Value = 9030
$ws.Range("I7").Value = 26277
$ws.Range("J7").Value = 29320
$ws.Range("K7").Value = 27269

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 236
$ws.Range("K4").Value = 100
$ws.Range("K6").Value = 197
$ws.Range("K7").Value = 810
$ws.Range("K8").Value = 1785
$ws.Range("K10").Value = 166
$ws.Range("K15").Value = 277
$ws.Range("K19").Value = 783
$ws.Range("K29").Value = 1505
$ws.Range("K33").Value = 1145
$ws.Range("K36").Value = 356
$ws.Range("K37").Value = 898
$ws.Range("K42").Value = 1014
$ws.Range("K44").Value = 220
$ws.Range("K48").Value = 343
$ws.Range("K49").Value = 154
$ws.Range("K51").Value = 352
$ws.Range("K53").Value = 346
$ws.Range("K54").Value = 530
$ws.Range("K58").Value = 17
$ws.Range("I63").Value = 239
$ws.Range("J63").Value = 197
$ws.Range("K63").Value = 78
$ws.Range("K65").Value = 632
$ws.Range("K67").Value = 1064
$ws.Range("K68").Value = 70
$ws.Range("K72").Value = 127
$ws.Range("K78").Value = 335
$ws.Range("K83").Value = 580
$ws.Range("K85").Value = 1260
$ws.Range("K89").Value = 409
$ws.Range("K91").Value = 329
$ws.Range("K93").Value = 110
$ws.Range("K99").Value = 462
$ws.Range("K100").Value = 48
$ws.Range("I101").Value = 26277
$ws.Range("J101").Value = 29320
$ws.Range("K101").Value = 27269

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 265
$ws.Range("K7").Value = 810

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K4").Value = 49
$ws.Range("K7").Value = 409

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 418
$ws.Range("K6").Value = 309
$ws.Range("K7").Value = 1260

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K6").Value = 144
$ws.Range("K7").Value = 346

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K3").Value = 541
$ws.Range("K5").Value = 48
$ws.Range("K6").Value = 600
$ws.Range("K7").Value = 1785

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K6").Value = 138
$ws.Range("K7").Value = 580

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K3").Value = 407
$ws.Range("K7").Value = 1145

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 257
$ws.Range("K3").Value = 299
$ws.Range("K7").Value = 898

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K6").Value = 235
$ws.Range("K7").Value = 632

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K3").Value = 196
$ws.Range("K7").Value = 462

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 291
$ws.Range("K7").Value = 1064

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("K6").Value = 78
$ws.Range("K7").Value = 154

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K2").Value = 83
$ws.Range("K6").Value = 285
$ws.Range("K7").Value = 530

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K3").Value = 531
$ws.Range("K7").Value = 1505

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K6").Value = 156
$ws.Range("K7").Value = 343

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 229
$ws.Range("K3").Value = 234
$ws.Range("K7").Value = 783

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K3").Value = 59
$ws.Range("K7").Value = 220

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K2").Value = 78
$ws.Range("K7").Value = 197

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K4").Value = 45
$ws.Range("K6").Value = 389
$ws.Range("K7").Value = 1014

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K2").Value = 52
$ws.Range("K7").Value = 166

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K6").Value = 107
$ws.Range("K7").Value = 335

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K3").Value = 150
$ws.Range("K4").Value = 20
$ws.Range("K7").Value = 329

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 217
$ws.Range("K4").Value = 41

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K2").Value = 134
$ws.Range("K7").Value = 356

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("K2").Value = 34
$ws.Range("K7").Value = 110

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("K6").Value = 28
$ws.Range("K7").Value = 48

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K2").Value = 104
$ws.Range("K7").Value = 277

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K3").Value = 69
$ws.Range("K7").Value = 236

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K6").Value = 114
$ws.Range("K7").Value = 352

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("K4").Value = 8
$ws.Range("K7").Value = 70

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K6").Value = 60
$ws.Range("K7").Value = 127

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("K2").Value = 34
$ws.Range("K3").Value = 22
$ws.Range("K7").Value = 100

$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Range("K3").Value = 2
$ws.Range("K7").Value = 17
